$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 26 appended after the existing data (A1:D25 -> A1:D26).
# Column A holds a literal date-like string ("2025-03-21"), not a real
# date value, matching the rest of the sheet (dates are stored as plain
# text). Excel's normal typed-value path auto-detects date-like strings
# and converts them to date serials, so we briefly force Text format,
# assign the literal string, then restore the default "Normal" style so
# the cell ends up with no explicit style (same as its neighbours) while
# keeping the value as text.
$dateCell = $ws.Cells.Item(26, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-03-21"
$dateCell.Style = "Normal"

$ws.Cells.Item(26, 2).Value = "Rien ne nous concerne aujourd'hui !"
$ws.Cells.Item(26, 3).Value = "NA"
$ws.Cells.Item(26, 4).Value = 1
